$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1: update the "datos actualizados" timestamp text
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 13 de Junio de 2020 a las 12:09"

# Row 46: 'Rumania' -> 'Rumania'
$ws.Cells.Item(46, 2).Value = 21679
$ws.Cells.Item(46, 3).Value = 275
$ws.Cells.Item(46, 4).Value = 15635
$ws.Cells.Item(46, 5).Value = 4650
$ws.Cells.Item(46, 7).Value = 14
$ws.Cells.Item(46, 8).Value = 1394

# Row 50: 'Barein' -> 'Barein'
$ws.Cells.Item(50, 5).Value = 5485
$ws.Cells.Item(50, 7).Value = 1
$ws.Cells.Item(50, 8).Value = 37

# Row 66: 'Camerun' -> 'Marruecos'
$ws.Cells.Item(66, 1).Value = "Marruecos"
$ws.Cells.Item(66, 2).Value = 8683
$ws.Cells.Item(66, 3).Value = 73
$ws.Cells.Item(66, 4).Value = 7664
$ws.Cells.Item(66, 5).Value = 807

# Row 67: 'Noruega' -> 'Camerun'
$ws.Cells.Item(67, 1).Value = "Camerun"
$ws.Cells.Item(67, 2).Value = 8681
$ws.Cells.Item(67, 4).Value = 4836
$ws.Cells.Item(67, 5).Value = 3633
$ws.Cells.Item(67, 8).Value = 212

# Row 68: 'Marruecos' -> 'Noruega'
$ws.Cells.Item(68, 1).Value = "Noruega"
$ws.Cells.Item(68, 2).Value = 8620
$ws.Cells.Item(68, 4).Value = 8138
$ws.Cells.Item(68, 5).Value = 240
$ws.Cells.Item(68, 8).Value = 242

# Row 72: 'Finlandia' -> 'Finlandia'
$ws.Cells.Item(72, 2).Value = 7087
$ws.Cells.Item(72, 3).Value = 14
$ws.Cells.Item(72, 5).Value = 562

# Row 78: 'Costa de Marfil' -> 'Consejo Danes para los Refugiados'
$ws.Cells.Item(78, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(78, 2).Value = 4724
$ws.Cells.Item(78, 3).Value = 87
$ws.Cells.Item(78, 4).Value = 595
$ws.Cells.Item(78, 5).Value = 4023
$ws.Cells.Item(78, 7).Value = 5
$ws.Cells.Item(78, 8).Value = 106

# Row 79: 'Consejo Danes para los Refugiados' -> 'Costa de Marfil'
$ws.Cells.Item(79, 1).Value = "Costa de Marfil"
$ws.Cells.Item(79, 2).Value = 4684
$ws.Cells.Item(79, 4).Value = 2263
$ws.Cells.Item(79, 5).Value = 2376
$ws.Cells.Item(79, 8).Value = 45

# Row 86: 'El Salvador' -> 'El Salvador'
$ws.Cells.Item(86, 2).Value = 3603
$ws.Cells.Item(86, 3).Value = 122
$ws.Cells.Item(86, 4).Value = 1738
$ws.Cells.Item(86, 5).Value = 1793

# Row 100: 'Republica de Africa Central' -> 'Republica de Africa Central'
$ws.Cells.Item(100, 2).Value = 2057
$ws.Cells.Item(100, 3).Value = 13
$ws.Cells.Item(100, 4).Value = 363
$ws.Cells.Item(100, 5).Value = 1687

# Row 112: 'Eslovenia' -> 'Eslovenia'
$ws.Cells.Item(112, 2).Value = 1492
$ws.Cells.Item(112, 3).Value = 2
$ws.Cells.Item(112, 5).Value = 24

# Row 114: 'Guinea-Bisau' -> 'Albania'
$ws.Cells.Item(114, 1).Value = "Albania"
$ws.Cells.Item(114, 2).Value = 1464
$ws.Cells.Item(114, 3).Value = 48
$ws.Cells.Item(114, 4).Value = 1039
$ws.Cells.Item(114, 5).Value = 389
$ws.Cells.Item(114, 8).Value = 36

# Row 115: 'Libano' -> 'Guinea-Bisau'
$ws.Cells.Item(115, 1).Value = "Guinea-Bisau"
$ws.Cells.Item(115, 2).Value = 1460
$ws.Cells.Item(115, 4).Value = 153
$ws.Cells.Item(115, 5).Value = 1292
$ws.Cells.Item(115, 8).Value = 15

# Row 116: 'Albania' -> 'Libano'
$ws.Cells.Item(116, 1).Value = "Libano"
$ws.Cells.Item(116, 2).Value = 1422
$ws.Cells.Item(116, 4).Value = 853
$ws.Cells.Item(116, 5).Value = 538
$ws.Cells.Item(116, 8).Value = 31

# Row 129: 'Burkina Faso' -> 'Burkina Faso'
$ws.Cells.Item(129, 4).Value = 799
$ws.Cells.Item(129, 5).Value = 40

# Row 206: 'Islas Malvinas' -> 'Groenlandia'
$ws.Cells.Item(206, 1).Value = "Groenlandia"

# Row 207: 'Groenlandia' -> 'Islas Malvinas'
$ws.Cells.Item(207, 1).Value = "Islas Malvinas"

# Row 208: 'Islas Turcas y Caicos' -> 'Santa Sede'
$ws.Cells.Item(208, 1).Value = "Santa Sede"
$ws.Cells.Item(208, 4).Value = 12
$ws.Cells.Item(208, 8).Value = 0

# Row 209: 'Santa Sede' -> 'Islas Turcas y Caicos'
$ws.Cells.Item(209, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(209, 4).Value = 11
$ws.Cells.Item(209, 8).Value = 1

# Row 210: 'Seychelles' -> 'Montserrat'
$ws.Cells.Item(210, 1).Value = "Montserrat"
$ws.Cells.Item(210, 4).Value = 10
$ws.Cells.Item(210, 8).Value = 1

# Row 211: 'Montserrat' -> 'Seychelles'
$ws.Cells.Item(211, 1).Value = "Seychelles"
$ws.Cells.Item(211, 4).Value = 11
$ws.Cells.Item(211, 8).Value = 0

# Row 213: 'Papua Nueva Guinea' -> 'Islas Virgenes Britanicas'
$ws.Cells.Item(213, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(213, 4).Value = 7
$ws.Cells.Item(213, 8).Value = 1

# Row 214: 'Islas Virgenes Britanicas' -> 'Papua Nueva Guinea'
$ws.Cells.Item(214, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(214, 4).Value = 8
$ws.Cells.Item(214, 8).Value = 0
